$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.7602039999999999
$arr[0,3] = 2.280612
$arr[0,4] = 0.9081302434927829
$arr[0,5] = 0.9140443705567521
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.887307
$arr[0,9] = 2.661921
$arr[0,10] = 0.03991909470044044
$arr[0,11] = 0.04024932703229714
$arr[0,12] = 0.6745343306279998
$arr[0,13] = 6.070808975651999
$arr[0,14] = 0.03625173719032244
$arr[0,15] = 0.03678967079256891
$ws.Range("E2:T2").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.7602039999999999
$arr[0,3] = 2.280612
$arr[0,4] = 0.9081302434927829
$arr[0,5] = 0.9140443705567521
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.9845453333333333
$arr[0,9] = 2.953636
$arr[0,10] = 0.04429375447078636
$arr[0,11] = 0.04466017635322986
$arr[0,12] = 0.7484553005813331
$arr[0,13] = 6.736097705231999
$arr[0,14] = 0.04022449803276476
$arr[0,15] = 0.04082138278374154
$ws.Range("E3:T3").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.7602039999999999
$arr[0,3] = 2.280612
$arr[0,4] = 0.9081302434927829
$arr[0,5] = 0.9140443705567521
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 15.43767133333333
$arr[0,9] = 46.313014
$arr[0,10] = 0.6945260928963797
$arr[0,11] = 0.7002715882016618
$arr[0,12] = 11.73577949828533
$arr[0,13] = 105.622015484568
$arr[0,14] = 0.6307201498540804
$arr[0,15] = 0.6400793030565651
$ws.Range("E4:T4").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.7602039999999999
$arr[0,3] = 2.280612
$arr[0,4] = 0.9081302434927829
$arr[0,5] = 0.9140443705567521
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 0.547111
$arr[0,9] = 1.094222
$arr[0,10] = 0.02461400149063703
$arr[0,11] = 0.01654508121162658
$arr[0,12] = 0.4159159706439999
$arr[0,13] = 2.495495823864
$arr[0,14] = 0.02235271916702393
$arr[0,15] = 0.01512293834189156
$ws.Range("E5:T5").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.7602039999999999
$arr[0,3] = 2.280612
$arr[0,4] = 0.9081302434927829
$arr[0,5] = 0.9140443705567521
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 4.370998666666666
$arr[0,9] = 13.112996
$arr[0,10] = 0.1966470564417564
$arr[0,11] = 0.1982738272011845
$arr[0,12] = 3.322850670394666
$arr[0,13] = 29.90565603355199
$arr[0,14] = 0.1785811392485913
$arr[0,15] = 0.1812310755819849
$ws.Range("E6:T6").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 1
$arr[0,1] = 0.5
$arr[0,2] = 0.016249
$arr[0,3] = 0.032498
$arr[0,4] = 0.01941085330584189
$arr[0,5] = 0.01302484331151171
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.887307
$arr[0,9] = 2.661921
$arr[0,10] = 0.03991909470044044
$arr[0,11] = 0.04024932703229714
$arr[0,12] = 0.014417851443
$arr[0,13] = 0.086507108658
$arr[0,14] = 0.00077486369133226
$arr[0,15] = 0.0005242411779894627
$ws.Range("E7:T7").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 1
$arr[0,1] = 0.5
$arr[0,2] = 0.016249
$arr[0,3] = 0.032498
$arr[0,4] = 0.01941085330584189
$arr[0,5] = 0.01302484331151171
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.9845453333333333
$arr[0,9] = 2.953636
$arr[0,10] = 0.04429375447078636
$arr[0,11] = 0.04466017635322986
$arr[0,12] = 0.01599787712133333
$arr[0,13] = 0.095987262728
$arr[0,14] = 0.0008597795703974125
$arr[0,15] = 0.0005816917992652993
$ws.Range("E8:T8").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 1
$arr[0,1] = 0.5
$arr[0,2] = 0.016249
$arr[0,3] = 0.032498
$arr[0,4] = 0.01941085330584189
$arr[0,5] = 0.01302484331151171
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 15.43767133333333
$arr[0,9] = 46.313014
$arr[0,10] = 0.6945260928963797
$arr[0,11] = 0.7002715882016618
$arr[0,12] = 0.2508467214953333
$arr[0,13] = 1.505080328972
$arr[0,14] = 0.01348134410629115
$arr[0,15] = 0.009120927711830094
$ws.Range("E9:T9").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 1
$arr[0,1] = 0.5
$arr[0,2] = 0.016249
$arr[0,3] = 0.032498
$arr[0,4] = 0.01941085330584189
$arr[0,5] = 0.01302484331151171
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 0.547111
$arr[0,9] = 1.094222
$arr[0,10] = 0.02461400149063703
$arr[0,11] = 0.01654508121162658
$arr[0,12] = 0.008890006639
$arr[0,13] = 0.035560026556
$arr[0,14] = 0.0004777787722045291
$arr[0,15] = 0.0002154970903576725
$ws.Range("E10:T10").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 1
$arr[0,1] = 0.5
$arr[0,2] = 0.016249
$arr[0,3] = 0.032498
$arr[0,4] = 0.01941085330584189
$arr[0,5] = 0.01302484331151171
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 4.370998666666666
$arr[0,9] = 13.112996
$arr[0,10] = 0.1966470564417564
$arr[0,11] = 0.1982738272011845
$arr[0,12] = 0.07102435733466665
$arr[0,13] = 0.4261461440079999
$arr[0,14] = 0.003817087165616544
$arr[0,15] = 0.002582485532069176
$ws.Range("E11:T11").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.06065600000000001
$arr[0,3] = 0.181968
$arr[0,4] = 0.07245890320137523
$arr[0,5] = 0.07293078613173619
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.887307
$arr[0,9] = 2.661921
$arr[0,10] = 0.03991909470044044
$arr[0,11] = 0.04024932703229714
$arr[0,12] = 0.05382049339200001
$arr[0,13] = 0.4843844405280001
$arr[0,14] = 0.002892493818785745
$arr[0,15] = 0.002935415061738771
$ws.Range("E12:T12").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.06065600000000001
$arr[0,3] = 0.181968
$arr[0,4] = 0.07245890320137523
$arr[0,5] = 0.07293078613173619
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.9845453333333333
$arr[0,9] = 2.953636
$arr[0,10] = 0.04429375447078636
$arr[0,11] = 0.04466017635322986
$arr[0,12] = 0.05971858173866667
$arr[0,13] = 0.537467235648
$arr[0,14] = 0.00320947686762419
$arr[0,15] = 0.003257101770223029
$ws.Range("E13:T13").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.06065600000000001
$arr[0,3] = 0.181968
$arr[0,4] = 0.07245890320137523
$arr[0,5] = 0.07293078613173619
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 15.43767133333333
$arr[0,9] = 46.313014
$arr[0,10] = 0.6945260928963797
$arr[0,11] = 0.7002715882016618
$arr[0,12] = 0.9363873923946667
$arr[0,13] = 8.427486531552
$arr[0,14] = 0.05032459893600812
$arr[0,15] = 0.05107135743326663
$ws.Range("E14:T14").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.06065600000000001
$arr[0,3] = 0.181968
$arr[0,4] = 0.07245890320137523
$arr[0,5] = 0.07293078613173619
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 0.547111
$arr[0,9] = 1.094222
$arr[0,10] = 0.02461400149063703
$arr[0,11] = 0.01654508121162658
$arr[0,12] = 0.033185564816
$arr[0,13] = 0.199113388896
$arr[0,14] = 0.001783503551408574
$arr[0,15] = 0.001206645779377345
$ws.Range("E15:T15").Value = $arr

$arr = New-Object "object[,]" 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.06065600000000001
$arr[0,3] = 0.181968
$arr[0,4] = 0.07245890320137523
$arr[0,5] = 0.07293078613173619
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 4.370998666666666
$arr[0,9] = 13.112996
$arr[0,10] = 0.1966470564417564
$arr[0,11] = 0.1982738272011845
$arr[0,12] = 0.2651272951253333
$arr[0,13] = 2.386145656128
$arr[0,14] = 0.0142488300275486
$arr[0,15] = 0.0144602660871304
$ws.Range("E16:T16").Value = $arr

